# The schedule table has a row whose "Заняття"/"Здача" date cells read
# 16.04 / 19.04 / 20.04, followed by an empty cell that (per the diff)
# should receive the text "25.04" in the same Times New Roman / 28half-pt
# formatting used throughout the table.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the target row by its known, unique cell contents instead of a
# hard-coded row number, so the script is resilient to minor layout shifts.
$targetRow = -1
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $c1 = $t.Cell($r, 1).Range.Text -replace "[\r\a]", ""
    $c2 = $t.Cell($r, 2).Range.Text -replace "[\r\a]", ""
    $c3 = $t.Cell($r, 3).Range.Text -replace "[\r\a]", ""
    $c4 = $t.Cell($r, 4).Range.Text -replace "[\r\a]", ""
    if ($c1 -eq "16.04" -and $c2 -eq "19.04" -and $c3 -eq "20.04" -and $c4 -eq "") {
        $targetRow = $r
        break
    }
}

if ($targetRow -eq -1) {
    throw "Could not locate the target schedule row (16.04 / 19.04 / 20.04 / <empty>)."
}

# Insert the new text into the empty cell.
$cell = $t.Cell($targetRow, 4)
$cell.Range.InsertAfter("25.04")

# Give the complex-script font name (w:cs) the same "Times New Roman" value
# as the rest of the table via a formatted Find/Replace pass scoped to the
# text we just inserted.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Replacement.Font.NameBi = "Times New Roman"
$find.Text = "25.04"
$find.Replacement.Text = ""
$find.Execute("25.04", $false, $false, $false, $false, $false, $true, 1, $true, "25.04", 2)

# Set the remaining (ascii/hAnsi + size + complex-script size) font
# attributes directly on the run so it matches sibling cells exactly:
# Times New Roman, 14pt (half-point value 28), including w:szCs.
$t.Cell($targetRow, 4).Range.Font.Name = "Times New Roman"
$t.Cell($targetRow, 4).Range.Font.Size = 14
$t.Cell($targetRow, 4).Range.Font.SizeBi = 14
